# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update the "K" column (G) values for each row of the save_data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 3
    3  = 3
    4  = 4
    5  = 2
    6  = 2
    7  = 3
    8  = 7
    9  = 1
    10 = 1
    11 = 1
    12 = 1
    13 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
